$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above row 10. This shifts:
#      old row 10 (P10:Q10 total)      -> row 11
#      old row 11 (footer / timestamp) -> row 12
$ws.Rows("10").Insert()

# 2) Build the new row 10 as a data row matching the layout/format of rows 7-9
#    by copying the formats from row 9 (values only, no merges) onto row 10,
#    then writing the new item's values and re-applying the merges.
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A10").Value2 = 4
$ws.Range("C10").Value2 = "اولويز ماكس طويل جدا"
$ws.Range("H10").Value2 = "37:0"

# L10 and P10 land on number-oriented formats (copied from L9/P9), so a
# plain numeric-looking string assignment gets silently coerced to a
# number. Force text storage, then restore the original number format so
# the cell's style/format is unaffected.
$origL10Fmt = $ws.Range("L10").NumberFormat
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value2 = "0"
$ws.Range("L10").NumberFormat = $origL10Fmt

$ws.Range("N10").Value2 = "35.00"

$origP10Fmt = $ws.Range("P10").NumberFormat
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value2 = "35.0000"
$ws.Range("P10").NumberFormat = $origP10Fmt

$ws.Range("Q10").Value2 = "1:0"

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Rows(10).RowHeight = 24.75

# 3) Update the (shifted) total row 11 value and height
$ws.Range("P11").Value2 = 176
$ws.Rows(11).RowHeight = 25.5

# 4) Update the (shifted) footer row 12 timestamp text
$ws.Range("A12").Value2 = "Wednesday, 8 October, 2025 9:38 AM"

"done"
